# Ausarbeitung/Paperprototype_Klassendiagramm.pptx — apply commit:
# "Hinzufügen von schriftlichem Teil der Ausarbeitung"
#
#  1. Update the fixed "datetimeFigureOut" footer date (master + all
#     slide layouts) from 22.01.2021 -> 28.02.2021.
#  2. Fill in the empty subtitle placeholder on slide 1 with the
#     matriculation number "Matrikel: 5920414".

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "22.01.2021") {
                $tr.Text = "28.02.2021"
                $tr.LanguageID = "de-DE"
            }
        }
    }
}

# --- Slide master ---
Update-DatePlaceholder $p.SlideMaster.Shapes

# --- Every slide layout under the master ---
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# --- Slide 1: add matriculation number to the (empty) subtitle ---
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.Name -eq "Untertitel 2") {
        $tr = $shp.TextFrame.TextRange
        $tr.Text = "Matrikel: 5920414"
        $tr.LanguageID = "de-DE"
    }
}
